# Updates crypto price/volume data per the commit diff (Tue May 21 14:13:33 UTC 2024).
# Rows 42/43 (Maker/Cosmos) and 44/45 (Bittensor/Arweave) are swapped as part of the update,
# plus the full set of D (Price) / E (Volume 1h) values are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.707.63'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +5.47%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.788.27'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +22.59%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '616.73'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +7.81%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.52'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.31%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.787.63'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +22.65%  '

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.546'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +6.38%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +12.67%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.39'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.59%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.503'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +7.69%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.79'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +13.58%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000263'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +9.39%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.407.24'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +22.22%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.793.19'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +22.66%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.852.18'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +5.70%  '

# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.65%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.61'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +8.65%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '525.02'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.03%  '

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.67%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.46'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +22.51%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.752'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +9.91%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.33'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.93%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.51'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +11.94%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.70'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +7.80%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.97'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +6.93%  '

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.04%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000127'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +35.38%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.51'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +9.80%  '

# Row 31
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +12.89%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.98'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.62%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.36'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +15.68%  '

# Row 34
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.94%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.14%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.21'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +11.80%  '

# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +11.44%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.345'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +11.35%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.21'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +10.03%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.133'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.18%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.72'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.72%  '

# Row 42
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.92'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +8.50%  '

# Row 43
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.153.01'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +12.56%  '

# Row 44
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '44.44'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.65%  '

# Row 45
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '426.10'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +15.36%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.83'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.06%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0371'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +8.02%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '28.02'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +8.96%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.76'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.29%  '

# Row 50
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.52'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +9.47%  '

# Row 51
$ws.Range('B51').Value = 'USDe'
$ws.Range('C51').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.00%  '
